$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 3.276109556179335
$ws.Range("C2").Value = 0.1710648090743803
$ws.Range("D2").Value = 0.05088992306788498
$ws.Range("E2").Value = 0.04304938155365168
$ws.Range("F2").Value = 6.946608703946595
$ws.Range("J2").Value = 0.1678877235439842
$ws.Range("L2").Value = 0.2796376590314011
$ws.Range("M2").Value = 0.6080991648005636
$ws.Range("B3").Value = 3.240473125576898
$ws.Range("C3").Value = 0.1593196973514921
$ws.Range("D3").Value = 0.04484173880580045
$ws.Range("E3").Value = 0.04265931115016475
$ws.Range("F3").Value = 6.776167955705859
$ws.Range("J3").Value = 0.16631733863656
$ws.Range("L3").Value = 0.2815785064640437
$ws.Range("M3").Value = 0.6051180434261099
$ws.Range("B4").Value = 3.221219722280807
$ws.Range("C4").Value = 0.1523128592110652
$ws.Range("D4").Value = 0.04112014866984737
$ws.Range("E4").Value = 0.04241428963117233
$ws.Range("F4").Value = 6.67300045888274
$ws.Range("J4").Value = 0.1653364572585012
$ws.Range("L4").Value = 0.2829349338563247
$ws.Range("M4").Value = 0.6037491444177547
$ws.Range("B5").Value = 3.214033277698718
$ws.Range("C5").Value = 0.1495084758910821
$ws.Range("D5").Value = 0.03960130796406247
$ws.Range("E5").Value = 0.04231303790313579
$ws.Range("F5").Value = 6.63132572387866
$ws.Range("J5").Value = 0.1649324567850741
$ws.Range("L5").Value = 0.2835291633455554
$ws.Range("M5").Value = 0.6033072436212734
$ws.Range("B6").Value = 3.212879769715528
$ws.Range("C6").Value = 0.1490458731736339
$ws.Range("D6").Value = 0.03934896152296119
$ws.Range("E6").Value = 0.04229613989232828
$ws.Range("F6").Value = 6.624427641420681
$ws.Range("J6").Value = 0.1648651113578694
$ws.Range("L6").Value = 0.2836303411312571
$ws.Range("M6").Value = 0.6032408653932535
$ws.Range("B7").Value = 3.221120134759701
$ws.Range("C7").Value = 0.1522748326543706
$ws.Range("D7").Value = 0.04109967454928665
$ws.Range("E7").Value = 0.04241292981896461
$ws.Range("F7").Value = 6.672436941662653
$ws.Range("J7").Value = 0.165331026234842
$ws.Range("L7").Value = 0.2829427798471116
$ws.Range("M7").Value = 0.6037427155093908
$ws.Range("B8").Value = 3.26327599087864
$ws.Range("C8").Value = 0.166972290154547
$ws.Range("D8").Value = 0.04880600971354454
$ws.Range("E8").Value = 0.04291602038801212
$ws.Range("F8").Value = 6.88752893112823
$ws.Range("J8").Value = 0.167349655161388
$ws.Range("L8").Value = 0.2802727002620244
$ws.Range("M8").Value = 0.6069754091130264
$ws.Range("B9").Value = 3.366860676288411
$ws.Range("C9").Value = 0.1974436504792152
$ws.Range("D9").Value = 0.06386695805305465
$ws.Range("E9").Value = 0.04385957043913713
$ws.Range("F9").Value = 7.321403551961964
$ws.Range("J9").Value = 0.1711803682579074
$ws.Range("L9").Value = 0.2763419084387948
$ws.Range("M9").Value = 0.6169838819496505
$ws.Range("B10").Value = 3.455831124156646
$ws.Range("C10").Value = 0.2208750269709014
$ws.Range("D10").Value = 0.07491882383149573
$ws.Range("E10").Value = 0.04452767904983634
$ws.Range("F10").Value = 7.647994882885484
$ws.Range("J10").Value = 0.1739230815192805
$ws.Range("L10").Value = 0.2742474856620802
$ws.Range("M10").Value = 0.6265867358435599
$ws.Range("B11").Value = 3.499127166959568
$ws.Range("C11").Value = 0.2317699354063336
$ws.Range("D11").Value = 0.07994760146419821
$ws.Range("E11").Value = 0.04482640719529041
$ws.Range("F11").Value = 7.798370922274728
$ws.Range("J11").Value = 0.1751565949933251
$ws.Range("L11").Value = 0.2734665935958063
$ws.Range("M11").Value = 0.631446927633931
$ws.Range("B12").Value = 3.51593020424508
$ws.Range("C12").Value = 0.2359301391719839
$ws.Range("D12").Value = 0.08185234598522584
$ws.Range("E12").Value = 0.04493879929569022
$ws.Range("F12").Value = 7.855582152492161
$ws.Range("J12").Value = 0.1756217652979934
$ws.Range("L12").Value = 0.273195573292071
$ws.Range("M12").Value = 0.6333583062042507
$ws.Range("B13").Value = 3.512293204034677
$ws.Range("C13").Value = 0.2350326186825953
$ws.Range("D13").Value = 0.08144210044009981
$ws.Range("E13").Value = 0.04491462581875361
$ws.Range("F13").Value = 7.843248697421018
$ws.Range("J13").Value = 0.1755216669564383
$ws.Range("L13").Value = 0.273252844844933
$ws.Range("M13").Value = 0.6329434985759974
$ws.Range("B14").Value = 3.500501379957143
$ws.Range("C14").Value = 0.2321115020550337
$ws.Range("D14").Value = 0.08010429520381024
$ws.Range("E14").Value = 0.04483566828458407
$ws.Range("F14").Value = 7.803072330737791
$ws.Range("J14").Value = 0.1751949030176512
$ws.Range("L14").Value = 0.2734438019899628
$ws.Range("M14").Value = 0.63160275513183
$ws.Range("B15").Value = 3.49333170600562
$ws.Range("C15").Value = 0.2303267504201756
$ws.Range("D15").Value = 0.07928491852551645
$ws.Range("E15").Value = 0.04478721000659291
$ws.Range("F15").Value = 7.778498141297007
$ws.Range("J15").Value = 0.1749945017807377
$ws.Range("L15").Value = 0.2735639828778176
$ws.Range("M15").Value = 0.6307907545529403
$ws.Range("B16").Value = 3.453058585346923
$ws.Range("C16").Value = 0.2201678180025226
$ws.Range("D16").Value = 0.07459022988494723
$ws.Range("E16").Value = 0.0445080533754636
$ws.Range("F16").Value = 7.638204441119484
$ws.Range("J16").Value = 0.1738421922478004
$ws.Range("L16").Value = 0.2743019748074786
$ws.Range("M16").Value = 0.626279024802642
$ws.Range("B17").Value = 3.429076630072075
$ws.Range("C17").Value = 0.2139964889545638
$ws.Range("D17").Value = 0.07171070266706181
$ws.Range("E17").Value = 0.0443354810447163
$ws.Range("F17").Value = 7.552606890770249
$ws.Range("J17").Value = 0.1731317329663362
$ws.Range("L17").Value = 0.2747987083806507
$ws.Range("M17").Value = 0.6236373323799214
$ws.Range("B18").Value = 3.415548437924315
$ws.Range("C18").Value = 0.2104690656429398
$ws.Range("D18").Value = 0.07005457355204214
$ws.Range("E18").Value = 0.04423573169113659
$ws.Range("F18").Value = 7.503543190871881
$ws.Range("J18").Value = 0.1727217565811934
$ws.Range("L18").Value = 0.2751005952801435
$ws.Range("M18").Value = 0.6221641786323246
$ws.Range("B19").Value = 3.411013586061131
$ws.Range("C19").Value = 0.2092785283329306
$ws.Range("D19").Value = 0.06949384520510193
$ws.Range("E19").Value = 0.04420187352742033
$ws.Range("F19").Value = 7.486960004764796
$ws.Range("J19").Value = 0.1725827125616632
$ws.Range("L19").Value = 0.2752055886385065
$ws.Range("M19").Value = 0.6216733361923303
$ws.Range("B20").Value = 3.431602049159096
$ws.Range("C20").Value = 0.2146511393306128
$ws.Range("D20").Value = 0.07201722102058739
$ws.Range("E20").Value = 0.04435390228209979
$ws.Range("F20").Value = 7.561701276874714
$ws.Range("J20").Value = 0.1732075004165026
$ws.Range("L20").Value = 0.2747441560040116
$ws.Range("M20").Value = 0.6239137539948132
$ws.Range("B21").Value = 3.503953842904195
$ws.Range("C21").Value = 0.2329685624219451
$ws.Range("D21").Value = 0.08049722642971346
$ws.Range("E21").Value = 0.04485887968991875
$ws.Range("F21").Value = 7.814865797725247
$ws.Range("J21").Value = 0.1752909331926418
$ws.Range("L21").Value = 0.2733870434740098
$ws.Range("M21").Value = 0.6319946372563194
$ws.Range("B22").Value = 3.553617598480514
$ws.Range("C22").Value = 0.2451416730313838
$ws.Range("D22").Value = 0.08604218711604972
$ws.Range("E22").Value = 0.04518466635449236
$ws.Range("F22").Value = 7.981883610621992
$ws.Range("J22").Value = 0.1766413492328418
$ws.Range("L22").Value = 0.2726439728954873
$ws.Range("M22").Value = 0.637689497882775
$ws.Range("B23").Value = 3.526892934420573
$ws.Range("C23").Value = 0.2386259974129246
$ws.Range("D23").Value = 0.08308238932477252
$ws.Range("E23").Value = 0.04501117042337555
$ws.Range("F23").Value = 7.892597931706575
$ws.Range("J23").Value = 0.1759216002455801
$ws.Range("L23").Value = 0.2730274072281205
$ws.Range("M23").Value = 0.6346121329801448
$ws.Range("B24").Value = 3.430459499499193
$ws.Range("C24").Value = 0.2143551079155657
$ws.Range("D24").Value = 0.07187864612011197
$ws.Range("E24").Value = 0.04434557570912823
$ws.Range("F24").Value = 7.557589246109785
$ws.Range("J24").Value = 0.1731732507063164
$ws.Range("L24").Value = 0.274768768334873
$ws.Range("M24").Value = 0.6237886417943486
$ws.Range("B25").Value = 3.336587318046725
$ws.Range("C25").Value = 0.1890198127253768
$ws.Range("D25").Value = 0.05979622062120882
$ws.Range("E25").Value = 0.04360880639190512
$ws.Range("F25").Value = 7.202688595208826
$ws.Range("J25").Value = 0.1701571032340006
$ws.Range("L25").Value = 0.2772657915718426
$ws.Range("M25").Value = 0.6138822324665156
